$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Tao%Guo%NULL%0,                Yongzhen%Fan%NULL%2,                Ming%Chen%NULL%3,                Xiaoyan%Wu%NULL%3,                Lin%Zhang%NULL%3,                Tao%He%NULL%3,                Hairong%Wang%NULL%3,                Jing%Wan%NULL%2,                Xinghuan%Wang%NULL%0,                Zhibing%Lu%NULL%3]'
$ws.Range("E3").Value = '[Suxin%Wan%NULL%0,                Yi%Xiang%NULL%0,                Wei%Fang%NULL%0,                Yu%Zheng%yuzheng1@cdutcm.edu.cn%0,                Boqun%Li%279685211@qq.com%0,                Boqun%Li%279685211@qq.com%0,                Yanjun%Hu%huyanjun@163.com%0,                Chunhui%Lang%NULL%0,                Daoqiu%Huang%NULL%0,                Qiuyan%Sun%NULL%0,                Yan%Xiong%NULL%0,                Xia%Huang%NULL%0,                Jinglong%Lv%NULL%0,                Yaling%Luo%NULL%0,                Li%Shen%NULL%0,                Haoran%Yang%NULL%0,                Gu%Huang%NULL%0,                Ruishan%Yang%NULL%0]'
$ws.Range("E4").Value = '[Fei%Zhou%NULL%0,                Ting%Yu%NULL%0,                Ronghui%Du%NULL%0,                Guohui%Fan%NULL%0,                Ying%Liu%NULL%0,                Zhibo%Liu%NULL%0,                Jie%Xiang%NULL%0,                Yeming%Wang%NULL%0,                Bin%Song%NULL%0,                Xiaoying%Gu%NULL%0,                Lulu%Guan%NULL%0,                Yuan%Wei%NULL%0,                Hui%Li%NULL%0,                Xudong%Wu%NULL%0,                Jiuyang%Xu%NULL%0,                Shengjin%Tu%NULL%0,                Yi%Zhang%NULL%0,                Hua%Chen%NULL%0,                Bin%Cao%NULL%0]'
$ws.Range("E6").Value = '[Yingzhen%Du%NULL%0,                Lei%Tu%NULL%5,                Lei%Tu%NULL%0,                Pingjun%Zhu%NULL%4,                Pingjun%Zhu%NULL%0,                Mi%Mu%NULL%5,                Mi%Mu%NULL%0,                Runsheng%Wang%NULL%3,                Pengcheng%Yang%NULL%3,                Xi%Wang%NULL%2,                Chao%Hu%NULL%3,                Rongyu%Ping%NULL%3,                Peng%Hu%NULL%3,                Tianzhi%Li%NULL%3,                Feng%Cao%NULL%2,                Christopher%Chang%NULL%2,                Qinyong%Hu%NULL%3,                Yang%Jin%NULL%2,                Guogang%Xu%NULL%3]'
$ws.Range("E7").Value = '[Liang%Chen%liang.chen9@hotmail.com%1,                Xiangjie%Li%NULL%1,                Mingquan%Chen%NULL%2,                Mingquan%Chen%NULL%0,                Yi%Feng%NULL%1,                Chenglong%Xiong%xiongchenglong@fudan.edu.cn%1]'
$ws.Range("E8").Value = '[Wei-jie%Guan%NULL%0,                Zheng-yi%Ni%NULL%0,                Zheng-yi%Ni%NULL%0,                Yu%Hu%NULL%0,                Wen-hua%Liang%NULL%0,                Chun-quan%Ou%NULL%0,                Jian-xing%He%NULL%0,                Lei%Liu%NULL%0,                Hong%Shan%NULL%0,                Chun-liang%Lei%NULL%0,                David S.C.%Hui%NULL%0,                Bin%Du%NULL%0,                Lan-juan%Li%NULL%0,                Guang%Zeng%NULL%0,                Kwok-Yung%Yuen%NULL%0,                Ru-chong%Chen%NULL%0,                Chun-li%Tang%NULL%0,                Tao%Wang%NULL%0,                Ping-yan%Chen%NULL%0,                Jie%Xiang%NULL%0,                Shi-yue%Li%NULL%0,                Jin-lin%Wang%NULL%0,                Zi-jing%Liang%NULL%0,                Yi-xiang%Peng%NULL%0,                Li%Wei%NULL%0,                Yong%Liu%NULL%0,                Ya-hua%Hu%NULL%0,                Peng%Peng%NULL%0,                Jian-ming%Wang%NULL%0,                Ji-yang%Liu%NULL%0,                Zhong%Chen%NULL%0,                Gang%Li%NULL%0,                Zhi-jian%Zheng%NULL%0,                Shao-qin%Qiu%NULL%0,                Jie%Luo%NULL%0,                Chang-jiang%Ye%NULL%0,                Shao-yong%Zhu%NULL%0,                Nan-shan%Zhong%NULL%0]'
$ws.Range("E9").Value = '[Akdis%Cezmi A%coreGivesNoEmail%0,               Cao%Yi-yuan%coreGivesNoEmail%0,               Dong%Xiang%coreGivesNoEmail%0,               Gao%Ya-dong%coreGivesNoEmail%0,               Yan%You-qin%coreGivesNoEmail%0,               Yang%Yi-bin%coreGivesNoEmail%0,               Yuan%Ya-dong%coreGivesNoEmail%0,               Zhang%Jin-jin%coreGivesNoEmail%0]'
$ws.Range("E10").Value = '[Lang%Wang%NULL%0,                Wenbo%He%NULL%3,                Xiaomei%Yu%NULL%3,                Dalong%Hu%NULL%3,                Mingwei%Bao%NULL%3,                Huafen%Liu%NULL%3,                Jiali%Zhou%NULL%3,                Hong%Jiang%NULL%3]'
$ws.Range("E11").Value = '[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%0,               Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%0,               Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%0,               Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%0,               Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%0]'
$ws.Range("E12").Value = '[Nanshan%Chen%NULL%0,                Min%Zhou%NULL%0,                Xuan%Dong%NULL%0,                Jieming%Qu%NULL%0,                Fengyun%Gong%NULL%0,                Yang%Han%NULL%0,                Yang%Qiu%NULL%0,                Jingli%Wang%NULL%0,                Ying%Liu%NULL%0,                Yuan%Wei%NULL%0,                Jia''an%Xia%NULL%0,                Ting%Yu%NULL%0,                Xinxin%Zhang%NULL%0,                Li%Zhang%NULL%0]'
$ws.Range("E13").Value = '[Xiaobo%Yang%NULL%0,                Yuan%Yu%NULL%0,                Jiqian%Xu%NULL%0,                Huaqing%Shu%NULL%0,                Jia''an%Xia%NULL%0,                Hong%Liu%NULL%0,                Yongran%Wu%NULL%0,                Lu%Zhang%NULL%0,                Zhui%Yu%NULL%0,                Minghao%Fang%NULL%0,                Ting%Yu%NULL%0,                Yaxin%Wang%NULL%0,                Shangwen%Pan%NULL%0,                Xiaojing%Zou%NULL%0,                Shiying%Yuan%NULL%0,                You%Shang%NULL%0]'
$ws.Range("E14").Value = '[Anesi%Adriano%coreGivesNoEmail%0,               Bettinardi%Alessandra%coreGivesNoEmail%0,               Bonetti%Graziella%coreGivesNoEmail%0,               Borrelli%Gianluca%coreGivesNoEmail%0,               Fiordalisi%Gianfranco%coreGivesNoEmail%0,               Lippi%Giuseppe%coreGivesNoEmail%0,               Manelli%Filippo%coreGivesNoEmail%0,               Marino%Antonio%coreGivesNoEmail%0,               Menolfi%Annamaria%coreGivesNoEmail%0,               Patroni%Andrea%coreGivesNoEmail%0,               Saggini%Sara%coreGivesNoEmail%0,               Volpi%Roberta%coreGivesNoEmail%0]'
$ws.Range("E16").Value = '[Huan%Han%NULL%0,                Linlin%Xie%NULL%1,                Rui%Liu%NULL%1,                Jie%Yang%NULL%2,                Fang%Liu%NULL%1,                Kailang%Wu%NULL%1,                Lang%Chen%NULL%1,                Wei%Hou%NULL%2,                Yong%Feng%yongfeng@whu.edu.cn%1,                Chengliang%Zhu%xinchengzhu@163.com%2,                Chengliang%Zhu%xinchengzhu@163.com%0]'
$ws.Range("E17").Value = '[Fan%Yang%NULL%2,                Shaobo%Shi%NULL%2,                Shaobo%Shi%NULL%0,                Jiling%Zhu%NULL%1,                Jinzhi%Shi%NULL%1,                Kai%Dai%NULL%1,                Xiaobei%Chen%121278395@qq.com%1]'
$ws.Range("E18").Value = '[Qing%Deng%NULL%2,                Bo%Hu%NULL%2,                Yao%Zhang%NULL%1,                Hao%Wang%NULL%1,                Xiaoyang%Zhou%NULL%1,                Wei%Hu%NULL%1,                Yuting%Cheng%NULL%1,                Jie%Yan%NULL%0,                Haiqin%Ping%NULL%1,                Qing%Zhou%qingzhou.wh.edu@hotmail.com%1]'
$ws.Range("E19").Value = '[Kyung Soo%Hong%NULL%2,                Kwan Ho%Lee%NULL%2,                Kwan Ho%Lee%NULL%0,                Jin Hong%Chung%NULL%2,                Jin Hong%Chung%NULL%0,                Kyeong-Cheol%Shin%NULL%2,                Kyeong-Cheol%Shin%NULL%0,                Eun Young%Choi%NULL%2,                Eun Young%Choi%NULL%0,                Hyun Jung%Jin%NULL%2,                Hyun Jung%Jin%NULL%0,                Jong Geol%Jang%NULL%2,                Jong Geol%Jang%NULL%0,                Wonhwa%Lee%NULL%2,                Wonhwa%Lee%NULL%0,                June Hong%Ahn%NULL%2,                June Hong%Ahn%NULL%0]'
$ws.Range("E20").Value = '[Yongli%Zheng%Zhyl1973@126.com%1,                Hong%Xu%NULL%1,                Ming%Yang%NULL%0,                Yilan%Zeng%NULL%1,                Hong%Chen%NULL%1,                Ru%Liu%NULL%1,                Qingfeng%Li%NULL%1,                Na%Zhang%NULL%1,                Dan%Wang%188170370@qq.com%2]'
